$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# Update Tuesday's Regular Hours from 1 to 2 (cell B12)
$ws.Range("B12").Value = 2

# Move the active selection to B13, matching the author's edit position
$ws.Range("B13").Select()
